$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 753.8461
$ws.Range("I8").Value = 66.666664
$ws.Range("J8").Value = 2300
$ws.Range("K8").Value = 199.999992
$ws.Range("L8").Value = 6900
$ws.Range("M8").Value = -60.99999199999999
$ws.Range("N8").Value = -7178

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1093.125
$ws.Range("I107").Value = 1307.8
$ws.Range("K107").Value = 1307.8
$ws.Range("M107").Value = 612.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3580.8
$ws.Range("I116").Value = 3580.8
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 3580.8
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -138.8000000000002
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 2250
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 26199.8
$ws.Range("I31").Value = 26199.8
$ws.Range("K31").Value = 26199.8
$ws.Range("M31").Value = -25905.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3070.6924
$ws.Range("I45").Value = 1756.2858
$ws.Range("K45").Value = 1756.2858
$ws.Range("M45").Value = -1379.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 756.8
$ws.Range("I97").Value = 729.7778
$ws.Range("K97").Value = 729.7778
$ws.Range("M97").Value = -233.7778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1441.7646
$ws.Range("I110").Value = 1441.8334
$ws.Range("J110").Value = 1441.6
$ws.Range("K110").Value = 1441.8334
$ws.Range("L110").Value = 1441.6
$ws.Range("M110").Value = 603.1666
$ws.Range("N110").Value = -5531.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2989
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2756
$ws.Range("I132").Value = 2692.7144
$ws.Range("K132").Value = 8078.1432
$ws.Range("M132").Value = -5548.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 31389
$ws.Range("I26").Value = 31389
$ws.Range("K26").Value = 31389
$ws.Range("M26").Value = -31097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 640
$ws.Range("I80").Value = 302.5
$ws.Range("J80").Value = 1765
$ws.Range("K80").Value = 302.5
$ws.Range("L80").Value = 1765
$ws.Range("M80").Value = 695.5
$ws.Range("N80").Value = -3761

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 640
$ws.Range("I83").Value = 302.5
$ws.Range("J83").Value = 1765
$ws.Range("K83").Value = 1512.5
$ws.Range("L83").Value = 8825
$ws.Range("M83").Value = 3479.5
$ws.Range("N83").Value = -18809

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 21084.834
$ws.Range("I96").Value = 21084.834
$ws.Range("K96").Value = 21084.834
$ws.Range("M96").Value = -18338.834

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1491.5834
$ws.Range("I105").Value = 1325.4445
$ws.Range("K105").Value = 1325.4445
$ws.Range("M105").Value = 421.5554999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 210000
$ws.Range("J138").Value = 210000
$ws.Range("L138").Value = 210000
$ws.Range("N138").Value = -220280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5703.4346
$ws.Range("I31").Value = 1813.75
$ws.Range("K31").Value = 1813.75
$ws.Range("M31").Value = -1518.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5703.4346
$ws.Range("I34").Value = 1813.75
$ws.Range("K34").Value = 1813.75
$ws.Range("M34").Value = -1611.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3719.7273
$ws.Range("I58").Value = 1920.5
$ws.Range("J58").Value = 5878.8
$ws.Range("K58").Value = 1920.5
$ws.Range("L58").Value = 5878.8
$ws.Range("M58").Value = -1717.5
$ws.Range("N58").Value = -6284.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1693.2858
$ws.Range("I132").Value = 1693.2858
$ws.Range("K132").Value = 5079.857400000001
$ws.Range("M132").Value = -2549.857400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1000
$ws.Range("I134").Value = 1000
$ws.Range("K134").Value = 3000
$ws.Range("M134").Value = -465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3719.7273
$ws.Range("I136").Value = 1920.5
$ws.Range("J136").Value = 5878.8
$ws.Range("K136").Value = 5761.5
$ws.Range("L136").Value = 17636.4
$ws.Range("M136").Value = -3211.5
$ws.Range("N136").Value = -22736.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 902
$ws.Range("J22").Value = 902
$ws.Range("L22").Value = 2706
$ws.Range("N22").Value = -3044

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 902
$ws.Range("J27").Value = 902
$ws.Range("L27").Value = 2706
$ws.Range("N27").Value = -2910

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1780
$ws.Range("I34").Value = 2614.5
$ws.Range("J34").Value = 1303.1428
$ws.Range("K34").Value = 7843.5
$ws.Range("L34").Value = 3909.4284
$ws.Range("M34").Value = -7759.5
$ws.Range("N34").Value = -4077.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 1300
$ws.Range("I49").Value = 100
$ws.Range("J49").Value = 2500
$ws.Range("K49").Value = 300
$ws.Range("L49").Value = 7500
$ws.Range("M49").Value = -144
$ws.Range("N49").Value = -7812

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 2286.25
$ws.Range("J82").Value = 2286.25
$ws.Range("L82").Value = 6858.75
$ws.Range("N82").Value = -7670.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H85").Value = 2286.25
$ws.Range("J85").Value = 2286.25
$ws.Range("L85").Value = 6858.75
$ws.Range("N85").Value = -9666.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1151.9
$ws.Range("J114").Value = 957.5
$ws.Range("L114").Value = 2872.5
$ws.Range("N114").Value = -9380.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1650
$ws.Range("I129").Value = 1650
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 4950
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 50
$ws.Range("N129").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2227.5715
$ws.Range("I102").Value = 2016.75
$ws.Range("K102").Value = 2016.75
$ws.Range("M102").Value = -394.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 75.75
$ws.Range("I107").Value = 75.75
$ws.Range("K107").Value = 75.75
$ws.Range("M107").Value = 1844.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 471.5
$ws.Range("I122").Value = 471.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1414.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 1035.5
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 8399.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5190.25
$ws.Range("I61").Value = 2436.25
$ws.Range("J61").Value = 7944.25
$ws.Range("K61").Value = 2436.25
$ws.Range("L61").Value = 7944.25
$ws.Range("M61").Value = -2234.25
$ws.Range("N61").Value = -8348.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 934.4706
$ws.Range("I93").Value = 947.7692
$ws.Range("J93").Value = 891.25
$ws.Range("K93").Value = 947.7692
$ws.Range("L93").Value = 891.25
$ws.Range("M93").Value = 300.2308
$ws.Range("N93").Value = -3387.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 5190.25
$ws.Range("I113").Value = 2436.25
$ws.Range("J113").Value = 7944.25
$ws.Range("K113").Value = 2436.25
$ws.Range("L113").Value = 7944.25
$ws.Range("M113").Value = -266.25
$ws.Range("N113").Value = -12284.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 776.05554
$ws.Range("I113").Value = 647.5
$ws.Range("K113").Value = 1942.5
$ws.Range("M113").Value = 227.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1333
$ws.Range("I122").Value = 1195.6
$ws.Range("J122").Value = 1447.5
$ws.Range("K122").Value = 3586.8
$ws.Range("L122").Value = 4342.5
$ws.Range("M122").Value = -1136.8
$ws.Range("N122").Value = -9242.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 979.5833
$ws.Range("I132").Value = 979.5833
$ws.Range("K132").Value = 2938.7499
$ws.Range("M132").Value = -408.7498999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1943.6072
$ws.Range("I136").Value = 1465.2
$ws.Range("K136").Value = 4395.6
$ws.Range("M136").Value = -1845.6
